$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "release/8.0.0"
$ws.Range("B3").Value = "X"
$ws.Range("C3").Value = "X"
$ws.Range("D3").Value = "X"
$ws.Range("E3").Value = "X"
